# Update "results with fixed workflow": shift the Cutoff/B column values by
# +4 (now starting at 5 instead of 1) and refresh the Reaction_number (C)
# column with newly computed results, then drop the last 4 rows (17-20) that
# no longer exist in the refreshed run, for both the NBR and BAR sheets.

$wb = $excel.ActiveWorkbook

$data = @{
    "NBR" = @(
        @(5,  883),
        @(6,  877),
        @(7,  879),
        @(8,  876),
        @(9,  877),
        @(10, 866),
        @(11, 862),
        @(12, 0),
        @(13, 909),
        @(14, 909),
        @(15, 905),
        @(16, 901),
        @(17, 887),
        @(18, 890),
        @(19, 886)
    )
    "BAR" = @(
        @(5,  854),
        @(6,  860),
        @(7,  856),
        @(8,  857),
        @(9,  839),
        @(10, 828),
        @(11, 830),
        @(12, 0),
        @(13, 757),
        @(14, 750),
        @(15, 743),
        @(16, 745),
        @(17, 745),
        @(18, 741),
        @(19, 742)
    )
}

foreach ($sheetName in @("NBR", "BAR")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 2).Value = $rows[$i][0]
        $ws.Cells.Item($r, 3).Value = $rows[$i][1]
    }

    $ws.Rows("17:20").Delete() | Out-Null
}
